$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 ("식물" theme): Fire count B8 increases 2 -> 3, and a new "ㅇ" mark
# is placed in the new "7-2" sub-position column (T8).
$ws.Range("B8").Value = 3
$ws.Range("T8").Value = "ㅇ"

# Row 9 ("식당" theme) is removed entirely - delete the whole row so the
# rows below (none, in this case) shift up and the used range shrinks to
# A1:V8.
$ws.Rows(9).Delete()

$ws.Range("T9").Select()
